$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2381892.8
$ws.Range("J17").Value = 2381892.8
$ws.Range("L17").Value = 7145678.399999999
$ws.Range("N17").Value = -7146014.399999999
$ws.Range("H42").Value = 525.0
$ws.Range("I42").Value = 1000.0
$ws.Range("J42").Value = 50.0
$ws.Range("K42").Value = 3000.0
$ws.Range("L42").Value = 150.0
$ws.Range("M42").Value = -2770.0
$ws.Range("N42").Value = -610.0
$ws.Range("H80").Value = 363.33334
$ws.Range("J80").Value = 296.0
$ws.Range("L80").Value = 888.0
$ws.Range("N80").Value = -2884.0
$ws.Range("H83").Value = 363.33334
$ws.Range("J83").Value = 296.0
$ws.Range("L83").Value = 2664.0
$ws.Range("N83").Value = -12648.0
$ws.Range("H118").Value = 995.0
$ws.Range("I118").Value = 995.0
$ws.Range("K118").Value = 2985.0
$ws.Range("M118").Value = -1328.0
$ws.Range("H131").Value = 130285.25
$ws.Range("I131").Value = 148040.28
$ws.Range("K131").Value = 444120.84
$ws.Range("M131").Value = -439080.84
$ws.Range("H137").Value = 1700.0
$ws.Range("I137").Value = 1050.0
$ws.Range("K137").Value = 3150.0
$ws.Range("M137").Value = -600.0
$ws.Range("H138").Value = 2366.4878
$ws.Range("J138").Value = 2505.4
$ws.Range("L138").Value = 7516.200000000001
$ws.Range("N138").Value = -17796.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2538.261
$ws.Range("I2").Value = 888.75
$ws.Range("K2").Value = 888.75
$ws.Range("M2").Value = -775.75
$ws.Range("H32").Value = 4021.9846
$ws.Range("I32").Value = 3928.5781
$ws.Range("J32").Value = 10000.0
$ws.Range("K32").Value = 3928.5781
$ws.Range("L32").Value = 10000.0
$ws.Range("M32").Value = -3641.5781
$ws.Range("N32").Value = -10574.0
$ws.Range("H45").Value = 6453.0
$ws.Range("I45").Value = 2975.0
$ws.Range("K45").Value = 2975.0
$ws.Range("M45").Value = -2598.0
$ws.Range("H74").Value = 3766.1035
$ws.Range("I74").Value = 3110.85
$ws.Range("J74").Value = 5222.222
$ws.Range("K74").Value = 3110.85
$ws.Range("L74").Value = 5222.222
$ws.Range("M74").Value = -2236.85
$ws.Range("N74").Value = -6970.222
$ws.Range("H77").Value = 3766.1035
$ws.Range("I77").Value = 3110.85
$ws.Range("J77").Value = 5222.222
$ws.Range("K77").Value = 15554.25
$ws.Range("L77").Value = 26111.11
$ws.Range("M77").Value = -11186.25
$ws.Range("N77").Value = -34847.11
$ws.Range("H97").Value = 784.26086
$ws.Range("I97").Value = 580.4286
$ws.Range("J97").Value = 2924.5
$ws.Range("K97").Value = 580.4286
$ws.Range("L97").Value = 2924.5
$ws.Range("M97").Value = -84.42859999999996
$ws.Range("N97").Value = -3916.5
$ws.Range("H116").Value = 2538.261
$ws.Range("I116").Value = 888.75
$ws.Range("K116").Value = 888.75
$ws.Range("M116").Value = 1405.25
$ws.Range("H132").Value = 5170.5835
$ws.Range("I132").Value = 4258.968
$ws.Range("K132").Value = 12776.904
$ws.Range("M132").Value = -10246.904

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2538.261
$ws.Range("I3").Value = 888.75
$ws.Range("K3").Value = 888.75
$ws.Range("M3").Value = -774.75
$ws.Range("H36").Value = 8996.0
$ws.Range("I36").Value = 8996.0
$ws.Range("J36").Value = 0.0
$ws.Range("K36").Value = 8996.0
$ws.Range("L36").Value = 0.0
$ws.Range("M36").Value = -8462.0
$ws.Range("H107").Value = 2086.7144
$ws.Range("I107").Value = 2086.7144
$ws.Range("K107").Value = 2086.7144
$ws.Range("M107").Value = -166.7143999999998
$ws.Range("H134").Value = 4521.375
$ws.Range("I134").Value = 4521.375
$ws.Range("K134").Value = 13564.125
$ws.Range("M134").Value = -11029.125
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5366.5806
$ws.Range("J31").Value = 5636.524
$ws.Range("L31").Value = 5636.524
$ws.Range("N31").Value = -6226.524
$ws.Range("H32").Value = 15000.0
$ws.Range("I32").Value = 15000.0
$ws.Range("K32").Value = 15000.0
$ws.Range("M32").Value = -14684.0
$ws.Range("H34").Value = 5366.5806
$ws.Range("J34").Value = 5636.524
$ws.Range("L34").Value = 5636.524
$ws.Range("N34").Value = -6040.524
$ws.Range("H86").Value = 10360.5
$ws.Range("I86").Value = 10721.5
$ws.Range("K86").Value = 10721.5
$ws.Range("M86").Value = -9598.5
$ws.Range("H89").Value = 10360.5
$ws.Range("I89").Value = 10721.5
$ws.Range("K89").Value = 53607.5
$ws.Range("M89").Value = -47991.5
$ws.Range("H122").Value = 3956.4119
$ws.Range("I122").Value = 3297.0715
$ws.Range("K122").Value = 9891.2145
$ws.Range("M122").Value = -7441.2145
$ws.Range("H133").Value = 52567.168
$ws.Range("I133").Value = 45000.0
$ws.Range("J133").Value = 54080.6
$ws.Range("K133").Value = 45000.0
$ws.Range("L133").Value = 54080.6
$ws.Range("M133").Value = -42470.0
$ws.Range("N133").Value = -59140.6
$ws.Range("H135").Value = 81388.0
$ws.Range("I135").Value = 0.0
$ws.Range("J135").Value = 81388.0
$ws.Range("K135").Value = 0.0
$ws.Range("L135").Value = 81388.0
$ws.Range("N135").Value = -91528.0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 143137.14
$ws.Range("J34").Value = 250375.0
$ws.Range("L34").Value = 751125.0
$ws.Range("N34").Value = -751293.0
$ws.Range("H40").Value = 119.181816
$ws.Range("I40").Value = 92.09091
$ws.Range("J40").Value = 146.27272
$ws.Range("K40").Value = 368.36364
$ws.Range("L40").Value = 585.09088
$ws.Range("M40").Value = -299.36364
$ws.Range("N40").Value = -723.09088
$ws.Range("H68").Value = 803.58826
$ws.Range("J68").Value = 779.8461
$ws.Range("L68").Value = 2339.5383
$ws.Range("N68").Value = -3961.5383
$ws.Range("H71").Value = 803.58826
$ws.Range("J71").Value = 779.8461
$ws.Range("L71").Value = 7018.6149
$ws.Range("N71").Value = -15130.6149
$ws.Range("H105").Value = 14995.0
$ws.Range("J105").Value = 14995.0
$ws.Range("L105").Value = 44985.0
$ws.Range("N105").Value = -50227.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4861.923
$ws.Range("I97").Value = 920.6
$ws.Range("K97").Value = 920.6
$ws.Range("M97").Value = -424.6
$ws.Range("H132").Value = 6832.4287
$ws.Range("J132").Value = 7467.1665
$ws.Range("L132").Value = 22401.4995
$ws.Range("N132").Value = -27461.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1535.6154
$ws.Range("J16").Value = 1713.4
$ws.Range("L16").Value = 1713.4
$ws.Range("N16").Value = -2053.4
$ws.Range("H22").Value = 1420.3889
$ws.Range("I22").Value = 1624.1111
$ws.Range("J22").Value = 1216.6666
$ws.Range("K22").Value = 1624.1111
$ws.Range("L22").Value = 1216.6666
$ws.Range("M22").Value = -1329.1111
$ws.Range("N22").Value = -1806.6666
$ws.Range("H27").Value = 1420.3889
$ws.Range("I27").Value = 1624.1111
$ws.Range("J27").Value = 1216.6666
$ws.Range("K27").Value = 1624.1111
$ws.Range("L27").Value = 1216.6666
$ws.Range("M27").Value = -1517.1111
$ws.Range("N27").Value = -1430.6666
$ws.Range("H29").Value = 7499.5
$ws.Range("I29").Value = 7499.5
$ws.Range("K29").Value = 7499.5
$ws.Range("M29").Value = -7204.5
$ws.Range("H33").Value = 25000.0
$ws.Range("I33").Value = 25000.0
$ws.Range("K33").Value = 25000.0
$ws.Range("M33").Value = -24710.0
$ws.Range("H40").Value = 3602.7
$ws.Range("I40").Value = 3256.8333
$ws.Range("J40").Value = 4121.5
$ws.Range("K40").Value = 3256.8333
$ws.Range("L40").Value = 4121.5
$ws.Range("M40").Value = -3120.8333
$ws.Range("N40").Value = -4393.5
$ws.Range("H61").Value = 54598.316
$ws.Range("I61").Value = 68705.664
$ws.Range("K61").Value = 68705.664
$ws.Range("M61").Value = -68503.664
$ws.Range("H113").Value = 54598.316
$ws.Range("I113").Value = 68705.664
$ws.Range("K113").Value = 68705.664
$ws.Range("M113").Value = -66535.664
$ws.Range("H122").Value = 3470.4849
$ws.Range("I122").Value = 3406.8
$ws.Range("K122").Value = 10220.4
$ws.Range("M122").Value = -7770.400000000001
$ws.Range("H137").Value = 0.0
$ws.Range("J137").Value = 0.0
$ws.Range("L137").Value = 0.0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 404.85715
$ws.Range("I107").Value = 222.33333
$ws.Range("K107").Value = 666.99999
$ws.Range("M107").Value = 1253.00001
$ws.Range("H124").Value = 0.0
$ws.Range("J124").Value = 0.0
$ws.Range("L124").Value = 0.0
$ws.Range("H132").Value = 2896.2954
$ws.Range("I132").Value = 2777.2368
$ws.Range("K132").Value = 8331.7104
$ws.Range("M132").Value = -5801.7104
$ws.Range("H136").Value = 4633.625
$ws.Range("I136").Value = 2665.55
$ws.Range("K136").Value = 8996.650000000001
$ws.Range("M136").Value = -5446.650000000001
$ws.Range("N124").ClearContents()
